$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Carry the existing formatting (bold+bordered "Indice" column style, and the
# datetime number format on "data_partida") down onto the two new rows by
# copying the formats from the last existing data row (86) before writing
# the new values.
$ws.Cells.Item(86, 1).Copy() | Out-Null
$ws.Cells.Item(87, 1).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(88, 1).PasteSpecial(-4122) | Out-Null

$ws.Cells.Item(86, 5).Copy() | Out-Null
$ws.Cells.Item(87, 5).PasteSpecial(-4122) | Out-Null
$ws.Cells.Item(88, 5).PasteSpecial(-4122) | Out-Null

$excel.CutCopyMode = 0

# Row 87: FC Copenhagen 0 x 0 Brondby
$ws.Cells.Item(87, 1).Value = 86
$ws.Cells.Item(87, 2).Value = "denmark"
$ws.Cells.Item(87, 3).Value = "superliga"
$ws.Cells.Item(87, 4).Value = "2023-2024"
$ws.Cells.Item(87, 5).Value = 45242.5
$ws.Cells.Item(87, 6).Value = "FC Copenhagen"
$ws.Cells.Item(87, 7).Value = 0
$ws.Cells.Item(87, 8).Value = "Brondby"
$ws.Cells.Item(87, 9).Value = 0
$ws.Cells.Item(87, 10).Value = 1.83
$ws.Cells.Item(87, 11).Value = "06/11/2023 19:12"
$ws.Cells.Item(87, 12).Value = 1.9
$ws.Cells.Item(87, 13).Value = "12/11/2023 11:59"
$ws.Cells.Item(87, 14).Value = 3.72
$ws.Cells.Item(87, 15).Value = "06/11/2023 19:12"
$ws.Cells.Item(87, 16).Value = 3.71
$ws.Cells.Item(87, 17).Value = "12/11/2023 11:59"
$ws.Cells.Item(87, 18).Value = 4
$ws.Cells.Item(87, 19).Value = "06/11/2023 19:12"
$ws.Cells.Item(87, 20).Value = 4.18
$ws.Cells.Item(87, 21).Value = "12/11/2023 11:59"
$ws.Cells.Item(87, 22).Value = "https://www.betexplorer.com/football/denmark/superliga/fc-copenhagen-brondby/2Lvp2ODd/"

# Row 88: Vejle 1 x 0 Lyngby
$ws.Cells.Item(88, 1).Value = 87
$ws.Cells.Item(88, 2).Value = "denmark"
$ws.Cells.Item(88, 3).Value = "superliga"
$ws.Cells.Item(88, 4).Value = "2023-2024"
$ws.Cells.Item(88, 5).Value = 45242.58333333334
$ws.Cells.Item(88, 6).Value = "Vejle"
$ws.Cells.Item(88, 7).Value = 1
$ws.Cells.Item(88, 8).Value = "Lyngby"
$ws.Cells.Item(88, 9).Value = 0
$ws.Cells.Item(88, 10).Value = 2.3
$ws.Cells.Item(88, 11).Value = "05/11/2023 14:12"
$ws.Cells.Item(88, 12).Value = 2.57
$ws.Cells.Item(88, 13).Value = "12/11/2023 13:52"
$ws.Cells.Item(88, 14).Value = 3.39
$ws.Cells.Item(88, 15).Value = "05/11/2023 14:12"
$ws.Cells.Item(88, 16).Value = 3.3
$ws.Cells.Item(88, 17).Value = "12/11/2023 13:52"
$ws.Cells.Item(88, 18).Value = 3.2
$ws.Cells.Item(88, 19).Value = "05/11/2023 14:12"
$ws.Cells.Item(88, 20).Value = 2.93
$ws.Cells.Item(88, 21).Value = "12/11/2023 13:52"
$ws.Cells.Item(88, 22).Value = "https://www.betexplorer.com/football/denmark/superliga/vejle-lyngby/fFRVPR53/"

Write-Output "Added rows 87-88"
